$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.0345785
$ws.Range("H2").Value = 0.069157
$ws.Range("M2").Value = 2.621797333333333
$ws.Range("N2").Value = 7.865392
$ws.Range("O2").Value = 0.07867217155043885
$ws.Range("P2").Value = 0.07906089226781998
$ws.Range("Q2").Value = 0.09065781909066666
$ws.Range("R2").Value = 0.543946914544
$ws.Range("S2").Value = 0.07867217155043885
$ws.Range("T2").Value = 0.07906089226781998

# Row 3
$ws.Range("G3").Value = 0.0345785
$ws.Range("H3").Value = 0.069157
$ws.Range("O3").Value = 0.148308476825081
$ws.Range("P3").Value = 0.1490412718702539
$ws.Range("Q3").Value = 0.1709031643165
$ws.Range("R3").Value = 1.025418985899
$ws.Range("S3").Value = 0.148308476825081
$ws.Range("T3").Value = 0.1490412718702539

# Row 4
$ws.Range("G4").Value = 0.0345785
$ws.Range("H4").Value = 0.069157
$ws.Range("M4").Value = 11.09754033333333
$ws.Range("N4").Value = 33.292621
$ws.Range("O4").Value = 0.3330034651388949
$ws.Range("P4").Value = 0.3346488416844782
$ws.Range("Q4").Value = 0.3837362984161666
$ws.Range("R4").Value = 2.302417790497
$ws.Range("S4").Value = 0.3330034651388949
$ws.Range("T4").Value = 0.3346488416844782

# Row 5
$ws.Range("G5").Value = 0.0345785
$ws.Range("H5").Value = 0.069157
$ws.Range("M5").Value = 0.4915585
$ws.Range("N5").Value = 0.983117
$ws.Range("O5").Value = 0.01475017696730553
$ws.Range("P5").Value = 0.009882038584175128
$ws.Range("Q5").Value = 0.01699735559225
$ws.Range("R5").Value = 0.06798942236899999
$ws.Range("S5").Value = 0.01475017696730553
$ws.Range("T5").Value = 0.009882038584175128

# Row 6
$ws.Range("G6").Value = 0.0345785
$ws.Range("H6").Value = 0.069157
$ws.Range("M6").Value = 14.172235
$ws.Range("N6").Value = 42.516705
$ws.Range("O6").Value = 0.4252657095182797
$ws.Range("P6").Value = 0.4273669555932728
$ws.Range("Q6").Value = 0.4900546279475
$ws.Range("R6").Value = 2.940327767685
$ws.Range("S6").Value = 0.4252657095182797
$ws.Range("T6").Value = 0.4273669555932728
